$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 0.1975515
$ws.Range("H2").Value = 0.395103
$ws.Range("I2").Value = 0.003060535523187668
$ws.Range("J2").Value = 0.002052987970621056
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 107.663086
$ws.Range("N2").Value = 215.326172
$ws.Range("O2").Value = 0.2751823527645522
$ws.Range("P2").Value = 0.2057131686791961
$ws.Range("Q2").Value = 21.269004133929
$ws.Range("R2").Value = 85.07601653571601
$ws.Range("S2").Value = 0.0008422053659902722
$ws.Range("T2").Value = 0.0004223266606967298

$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 0.1975515
$ws.Range("H3").Value = 0.395103
$ws.Range("I3").Value = 0.003060535523187668
$ws.Range("J3").Value = 0.002052987970621056
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.41172733333334
$ws.Range("N3").Value = 82.23518200000001
$ws.Range("O3").Value = 0.07006323059444025
$ws.Range("P3").Value = 0.07856388152449204
$ws.Range("Q3").Value = 5.415227852291001
$ws.Range("R3").Value = 32.49136711374601
$ws.Range("S3").Value = 0.0002144310061035734
$ws.Range("T3").Value = 0.00016129070369508

$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 0.1975515
$ws.Range("H4").Value = 0.395103
$ws.Range("I4").Value = 0.003060535523187668
$ws.Range("J4").Value = 0.002052987970621056
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 82.303927
$ws.Range("N4").Value = 246.911781
$ws.Range("O4").Value = 0.2103654011331419
$ws.Range("P4").Value = 0.2358886724356653
$ws.Range("Q4").Value = 16.2592642347405
$ws.Range("R4").Value = 97.55558540844301
$ws.Range("S4").Value = 0.0006438307830176041
$ws.Range("T4").Value = 0.0004842766069161916

$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 0.1975515
$ws.Range("H5").Value = 0.395103
$ws.Range("I5").Value = 0.003060535523187668
$ws.Range("J5").Value = 0.002052987970621056
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 71.272429
$ws.Range("N5").Value = 213.817287
$ws.Range("O5").Value = 0.1821693528222338
$ws.Range("P5").Value = 0.2042716462128862
$ws.Range("Q5").Value = 14.0799752575935
$ws.Range("R5").Value = 84.479851545561
$ws.Range("S5").Value = 0.0005575357755485543
$ws.Range("T5").Value = 0.0004193672324140156

$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 0.1975515
$ws.Range("H6").Value = 0.395103
$ws.Range("I6").Value = 0.003060535523187668
$ws.Range("J6").Value = 0.002052987970621056
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 83.25665266666667
$ws.Range("N6").Value = 249.769958
$ws.Range("O6").Value = 0.2128005281598046
$ws.Range("P6").Value = 0.2386192492246123
$ws.Range("Q6").Value = 16.447476619279
$ws.Range("R6").Value = 98.68485971567401
$ws.Range("S6").Value = 0.0006512835757861796
$ws.Range("T6").Value = 0.0004898824482167568

$ws.Range("E7").Value = 2
$ws.Range("G7").Value = 0.1975515
$ws.Range("H7").Value = 0.395103
$ws.Range("I7").Value = 0.003060535523187668
$ws.Range("J7").Value = 0.002052987970621056
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 19.3348755
$ws.Range("N7").Value = 38.66975100000001
$ws.Range("O7").Value = 0.04941913452582716
$ws.Range("P7").Value = 0.03694338192314826
$ws.Range("Q7").Value = 3.819633657338251
$ws.Range("R7").Value = 15.278534629353
$ws.Range("S7").Value = 0.0001512490167414842
$ws.Range("T7").Value = 0.00007584431868228277

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.293588333333333
$ws.Range("H8").Value = 3.880765
$ws.Range("I8").Value = 0.02004071366984202
$ws.Range("J8").Value = 0.02016477693615898
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 107.663086
$ws.Range("N8").Value = 215.326172
$ws.Range("O8").Value = 0.2751823527645522
$ws.Range("P8").Value = 0.2057131686791961
$ws.Range("Q8").Value = 139.2717119802633
$ws.Range("R8").Value = 835.63027188158
$ws.Range("S8").Value = 0.005514850738747851
$ws.Range("T8").Value = 0.004148160159246436

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.293588333333333
$ws.Range("H9").Value = 3.880765
$ws.Range("I9").Value = 0.02004071366984202
$ws.Range("J9").Value = 0.02016477693615898
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 27.41172733333334
$ws.Range("N9").Value = 82.23518200000001
$ws.Range("O9").Value = 0.07006323059444025
$ws.Range("P9").Value = 0.07856388152449204
$ws.Range("Q9").Value = 35.45949067491445
$ws.Range("R9").Value = 319.1354160742301
$ws.Range("S9").Value = 0.001404117143127292
$ws.Range("T9").Value = 0.001584223146180204

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.293588333333333
$ws.Range("H10").Value = 3.880765
$ws.Range("I10").Value = 0.02004071366984202
$ws.Range("J10").Value = 0.02016477693615898
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 82.303927
$ws.Range("N10").Value = 246.911781
$ws.Range("O10").Value = 0.2103654011331419
$ws.Range("P10").Value = 0.2358886724356653
$ws.Range("Q10").Value = 106.4673997547183
$ws.Range("R10").Value = 958.206597792465
$ws.Range("S10").Value = 0.004215872770150757
$ws.Range("T10").Value = 0.004756642461431865

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.293588333333333
$ws.Range("H11").Value = 3.880765
$ws.Range("I11").Value = 0.02004071366984202
$ws.Range("J11").Value = 0.02016477693615898
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 71.272429
$ws.Range("N11").Value = 213.817287
$ws.Range("O11").Value = 0.1821693528222338
$ws.Range("P11").Value = 0.2042716462128862
$ws.Range("Q11").Value = 92.19718264272834
$ws.Range("R11").Value = 829.774643784555
$ws.Range("S11").Value = 0.003650803839330815
$ws.Range("T11").Value = 0.004119092180264836

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.293588333333333
$ws.Range("H12").Value = 3.880765
$ws.Range("I12").Value = 0.02004071366984202
$ws.Range("J12").Value = 0.02016477693615898
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 83.25665266666667
$ws.Range("N12").Value = 249.769958
$ws.Range("O12").Value = 0.2128005281598046
$ws.Range("P12").Value = 0.2386192492246123
$ws.Range("Q12").Value = 107.6998345619856
$ws.Range("R12").Value = 969.2985110578701
$ws.Range("S12").Value = 0.004264674453641798
$ws.Range("T12").Value = 0.004811703933288034

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.293588333333333
$ws.Range("H13").Value = 3.880765
$ws.Range("I13").Value = 0.02004071366984202
$ws.Range("J13").Value = 0.02016477693615898
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 19.3348755
$ws.Range("N13").Value = 38.66975100000001
$ws.Range("O13").Value = 0.04941913452582716
$ws.Range("P13").Value = 0.03694338192314826
$ws.Range("Q13").Value = 25.0113693732525
$ws.Range("R13").Value = 150.068216239515
$ws.Range("S13").Value = 0.0009903947248435061
$ws.Range("T13").Value = 0.0007449550557476127

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 39.745772
$ws.Range("H14").Value = 119.237316
$ws.Range("I14").Value = 0.6157551175390606
$ws.Range("J14").Value = 0.6195669873353065
$ws.Range("K14").Value = 2
$ws.Range("M14").Value = 107.663086
$ws.Range("N14").Value = 215.326172
$ws.Range("O14").Value = 0.2751823527645522
$ws.Range("P14").Value = 0.2057131686791961
$ws.Range("Q14").Value = 4279.152468972392
$ws.Range("R14").Value = 25674.91481383435
$ws.Range("S14").Value = 0.1694449419712121
$ws.Range("T14").Value = 0.1274530881737692

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 39.745772
$ws.Range("H15").Value = 119.237316
$ws.Range("I15").Value = 0.6157551175390606
$ws.Range("J15").Value = 0.6195669873353065
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 27.41172733333334
$ws.Range("N15").Value = 82.23518200000001
$ws.Range("O15").Value = 0.07006323059444025
$ws.Range("P15").Value = 0.07856388152449204
$ws.Range("Q15").Value = 1089.500264716835
$ws.Range("R15").Value = 9805.502382451514
$ws.Range("S15").Value = 0.04314179278984587
$ws.Range("T15").Value = 0.04867558738949748

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 39.745772
$ws.Range("H16").Value = 119.237316
$ws.Range("I16").Value = 0.6157551175390606
$ws.Range("J16").Value = 0.6195669873353065
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 82.303927
$ws.Range("N16").Value = 246.911781
$ws.Range("O16").Value = 0.2103654011331419
$ws.Range("P16").Value = 0.2358886724356653
$ws.Range("Q16").Value = 3271.233117246644
$ws.Range("R16").Value = 29441.0980552198
$ws.Range("S16").Value = 0.1295335723008894
$ws.Range("T16").Value = 0.1461488341274901

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 39.745772
$ws.Range("H17").Value = 119.237316
$ws.Range("I17").Value = 0.6157551175390606
$ws.Range("J17").Value = 0.6195669873353065
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 71.272429
$ws.Range("N17").Value = 213.817287
$ws.Range("O17").Value = 0.1821693528222338
$ws.Range("P17").Value = 0.2042716462128862
$ws.Range("Q17").Value = 2832.777712920188
$ws.Range("R17").Value = 25494.99941628169
$ws.Range("S17").Value = 0.1121717112590692
$ws.Range("T17").Value = 0.1265599684421415

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 39.745772
$ws.Range("H18").Value = 119.237316
$ws.Range("I18").Value = 0.6157551175390606
$ws.Range("J18").Value = 0.6195669873353065
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 83.25665266666667
$ws.Range("N18").Value = 249.769958
$ws.Range("O18").Value = 0.2128005281598046
$ws.Range("P18").Value = 0.2386192492246123
$ws.Range("Q18").Value = 3309.099934372526
$ws.Range("R18").Value = 29781.89940935273
$ws.Range("S18").Value = 0.1310330142294147
$ws.Range("T18").Value = 0.1478406093623057

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 39.745772
$ws.Range("H19").Value = 119.237316
$ws.Range("I19").Value = 0.6157551175390606
$ws.Range("J19").Value = 0.6195669873353065
$ws.Range("K19").Value = 2
$ws.Range("M19").Value = 19.3348755
$ws.Range("N19").Value = 38.66975100000001
$ws.Range("O19").Value = 0.04941913452582716
$ws.Range("P19").Value = 0.03694338192314826
$ws.Range("Q19").Value = 768.4795532713862
$ws.Range("R19").Value = 4610.877319628317
$ws.Range("S19").Value = 0.03043008498862936
$ws.Range("T19").Value = 0.02288889984010259

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 22.079986
$ws.Range("H20").Value = 66.239958
$ws.Range("I20").Value = 0.3420707081671684
$ws.Range("J20").Value = 0.3441883178524182
$ws.Range("K20").Value = 2
$ws.Range("M20").Value = 107.663086
$ws.Range("N20").Value = 215.326172
$ws.Range("O20").Value = 0.2751823527645522
$ws.Range("P20").Value = 0.2057131686791961
$ws.Range("Q20").Value = 2377.199431596796
$ws.Range("R20").Value = 14263.19658958078
$ws.Range("S20").Value = 0.09413182228527794
$ws.Range("T20").Value = 0.07080406948778326

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 22.079986
$ws.Range("H21").Value = 66.239958
$ws.Range("I21").Value = 0.3420707081671684
$ws.Range("J21").Value = 0.3441883178524182
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 27.41172733333334
$ws.Range("N21").Value = 82.23518200000001
$ws.Range("O21").Value = 0.07006323059444025
$ws.Range("P21").Value = 0.07856388152449204
$ws.Range("Q21").Value = 605.2505557558175
$ws.Range("R21").Value = 5447.255001802357
$ws.Range("S21").Value = 0.0239665789059198
$ws.Range("T21").Value = 0.02704077022587159

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 22.079986
$ws.Range("H22").Value = 66.239958
$ws.Range("I22").Value = 0.3420707081671684
$ws.Range("J22").Value = 0.3441883178524182
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 82.303927
$ws.Range("N22").Value = 246.911781
$ws.Range("O22").Value = 0.2103654011331419
$ws.Range("P22").Value = 0.2358886724356653
$ws.Range("Q22").Value = 1817.269555905022
$ws.Range("R22").Value = 16355.4260031452
$ws.Range("S22").Value = 0.07195984173948432
$ws.Range("T22").Value = 0.08119012536607172

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 22.079986
$ws.Range("H23").Value = 66.239958
$ws.Range("I23").Value = 0.3420707081671684
$ws.Range("J23").Value = 0.3441883178524182
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 71.272429
$ws.Range("N23").Value = 213.817287
$ws.Range("O23").Value = 0.1821693528222338
$ws.Range("P23").Value = 0.2042716462128862
$ws.Range("Q23").Value = 1573.694234505994
$ws.Range("R23").Value = 14163.24811055395
$ws.Range("S23").Value = 0.06231479952625629
$ws.Range("T23").Value = 0.07030791429495759

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 22.079986
$ws.Range("H24").Value = 66.239958
$ws.Range("I24").Value = 0.3420707081671684
$ws.Range("J24").Value = 0.3441883178524182
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 83.25665266666667
$ws.Range("N24").Value = 249.769958
$ws.Range("O24").Value = 0.2128005281598046
$ws.Range("P24").Value = 0.2386192492246123
$ws.Range("Q24").Value = 1838.305725286863
$ws.Range("R24").Value = 16544.75152758176
$ws.Range("S24").Value = 0.07279282736597183
$ws.Range("T24").Value = 0.08212995799782624

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 22.079986
$ws.Range("H25").Value = 66.239958
$ws.Range("I25").Value = 0.3420707081671684
$ws.Range("J25").Value = 0.3441883178524182
$ws.Range("K25").Value = 2
$ws.Range("M25").Value = 19.3348755
$ws.Range("N25").Value = 38.66975100000001
$ws.Range("O25").Value = 0.04941913452582716
$ws.Range("P25").Value = 0.03694338192314826
$ws.Range("Q25").Value = 426.9137803517431
$ws.Range("R25").Value = 2561.482682110458
$ws.Range("S25").Value = 0.01690483834425826
$ws.Range("T25").Value = 0.01271548047990783

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 0.237281
$ws.Range("H26").Value = 0.711843
$ws.Range("I26").Value = 0.003676038549327608
$ws.Range("J26").Value = 0.003698795291280512
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 107.663086
$ws.Range("N26").Value = 215.326172
$ws.Range("O26").Value = 0.2751823527645522
$ws.Range("P26").Value = 0.2057131686791961
$ws.Range("Q26").Value = 25.546404709166
$ws.Range("R26").Value = 153.278428254996
$ws.Range("S26").Value = 0.001011580936857163
$ws.Range("T26").Value = 0.0007608908996650042

$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 0.237281
$ws.Range("H27").Value = 0.711843
$ws.Range("I27").Value = 0.003676038549327608
$ws.Range("J27").Value = 0.003698795291280512
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 27.41172733333334
$ws.Range("N27").Value = 82.23518200000001
$ws.Range("O27").Value = 0.07006323059444025
$ws.Range("P27").Value = 0.07856388152449204
$ws.Range("Q27").Value = 6.504282073380667
$ws.Range("R27").Value = 58.53853866042601
$ws.Range("S27").Value = 0.0002575551365555918
$ws.Range("T27").Value = 0.0002905917150475111

$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 0.237281
$ws.Range("H28").Value = 0.711843
$ws.Range("I28").Value = 0.003676038549327608
$ws.Range("J28").Value = 0.003698795291280512
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 82.303927
$ws.Range("N28").Value = 246.911781
$ws.Range("O28").Value = 0.2103654011331419
$ws.Range("P28").Value = 0.2358886724356653
$ws.Range("Q28").Value = 19.529158102487
$ws.Range("R28").Value = 175.762422922383
$ws.Range("S28").Value = 0.0007733113240101954
$ws.Range("T28").Value = 0.0008725039108714499

$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 0.237281
$ws.Range("H29").Value = 0.711843
$ws.Range("I29").Value = 0.003676038549327608
$ws.Range("J29").Value = 0.003698795291280512
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 71.272429
$ws.Range("N29").Value = 213.817287
$ws.Range("O29").Value = 0.1821693528222338
$ws.Range("P29").Value = 0.2042716462128862
$ws.Range("Q29").Value = 16.911593225549
$ws.Range("R29").Value = 152.204339029941
$ws.Range("S29").Value = 0.0006696615634805936
$ws.Range("T29").Value = 0.0007555590031543422

$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 0.237281
$ws.Range("H30").Value = 0.711843
$ws.Range("I30").Value = 0.003676038549327608
$ws.Range("J30").Value = 0.003698795291280512
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 83.25665266666667
$ws.Range("N30").Value = 249.769958
$ws.Range("O30").Value = 0.2128005281598046
$ws.Range("P30").Value = 0.2386192492246123
$ws.Range("Q30").Value = 19.75522180139933
$ws.Range("R30").Value = 177.796996212594
$ws.Range("S30").Value = 0.0007822629448327169
$ws.Range("T30").Value = 0.000882603755440887

$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 0.237281
$ws.Range("H31").Value = 0.711843
$ws.Range("I31").Value = 0.003676038549327608
$ws.Range("J31").Value = 0.003698795291280512
$ws.Range("K31").Value = 2
$ws.Range("M31").Value = 19.3348755
$ws.Range("N31").Value = 38.66975100000001
$ws.Range("O31").Value = 0.04941913452582716
$ws.Range("P31").Value = 0.03694338192314826
$ws.Range("Q31").Value = 4.5877985935155
$ws.Range("R31").Value = 27.526791561093
$ws.Range("S31").Value = 0.0001816666435913476
$ws.Range("T31").Value = 0.0001366460071013184

$ws.Range("E32").Value = 2
$ws.Range("G32").Value = 0.9938385
$ws.Range("H32").Value = 1.987677
$ws.Range("I32").Value = 0.01539688655141341
$ws.Range("J32").Value = 0.0103281346142149
$ws.Range("K32").Value = 2
$ws.Range("M32").Value = 107.663086
$ws.Range("N32").Value = 215.326172
$ws.Range("O32").Value = 0.2751823527645522
$ws.Range("P32").Value = 0.2057131686791961
$ws.Range("Q32").Value = 106.999719895611
$ws.Range("R32").Value = 427.998879582444
$ws.Range("S32").Value = 0.004236951466466835
$ws.Range("T32").Value = 0.002124633298035433

$ws.Range("E33").Value = 2
$ws.Range("G33").Value = 0.9938385
$ws.Range("H33").Value = 1.987677
$ws.Range("I33").Value = 0.01539688655141341
$ws.Range("J33").Value = 0.0103281346142149
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 27.41172733333334
$ws.Range("N33").Value = 82.23518200000001
$ws.Range("O33").Value = 0.07006323059444025
$ws.Range("P33").Value = 0.07856388152449204
$ws.Range("Q33").Value = 27.242829975369
$ws.Range("R33").Value = 163.456979852214
$ws.Range("S33").Value = 0.001078755612888114
$ws.Range("T33").Value = 0.0008114183442001844

$ws.Range("E34").Value = 2
$ws.Range("G34").Value = 0.9938385
$ws.Range("H34").Value = 1.987677
$ws.Range("I34").Value = 0.01539688655141341
$ws.Range("J34").Value = 0.0103281346142149
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 82.303927
$ws.Range("N34").Value = 246.911781
$ws.Range("O34").Value = 0.2103654011331419
$ws.Range("P34").Value = 0.2358886724356653
$ws.Range("Q34").Value = 81.79681135378949
$ws.Range("R34").Value = 490.7808681227369
$ws.Range("S34").Value = 0.00323897221558956
$ws.Range("T34").Value = 0.002436289962883994

$ws.Range("E35").Value = 2
$ws.Range("G35").Value = 0.9938385
$ws.Range("H35").Value = 1.987677
$ws.Range("I35").Value = 0.01539688655141341
$ws.Range("J35").Value = 0.0103281346142149
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 71.272429
$ws.Range("N35").Value = 213.817287
$ws.Range("O35").Value = 0.1821693528222338
$ws.Range("P35").Value = 0.2042716462128862
$ws.Range("Q35").Value = 70.8332839287165
$ws.Range("R35").Value = 424.999703572299
$ws.Range("S35").Value = 0.002804840858548336
$ws.Range("T35").Value = 0.002109745059953969

$ws.Range("E36").Value = 2
$ws.Range("G36").Value = 0.9938385
$ws.Range("H36").Value = 1.987677
$ws.Range("I36").Value = 0.01539688655141341
$ws.Range("J36").Value = 0.0103281346142149
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 83.25665266666667
$ws.Range("N36").Value = 249.769958
$ws.Range("O36").Value = 0.2128005281598046
$ws.Range("P36").Value = 0.2386192492246123
$ws.Range("Q36").Value = 82.743666801261
$ws.Range("R36").Value = 496.462000807566
$ws.Range("S36").Value = 0.003276465590157366
$ws.Range("T36").Value = 0.002464491727534689

$ws.Range("E37").Value = 2
$ws.Range("G37").Value = 0.9938385
$ws.Range("H37").Value = 1.987677
$ws.Range("I37").Value = 0.01539688655141341
$ws.Range("J37").Value = 0.0103281346142149
$ws.Range("K37").Value = 2
$ws.Range("M37").Value = 19.3348755
$ws.Range("N37").Value = 38.66975100000001
$ws.Range("O37").Value = 0.04941913452582716
$ws.Range("P37").Value = 0.03694338192314826
$ws.Range("Q37").Value = 19.21574366460675
$ws.Range("R37").Value = 76.86297465842701
$ws.Range("S37").Value = 0.0007609008077631984
$ws.Range("T37").Value = 0.0003815562216066285
